# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 156 (pushing existing rows 156-177
# down to 157-178) on the single data sheet, and populate it with the
# new week's observation for "Haba" at Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 156; everything below shifts down.
$ws.Range("A156:R156").EntireRow.Insert()

$ws.Cells.Item(156, 1).Value = 9
$ws.Cells.Item(156, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(156, 3).Value = "Metropolitana"
$ws.Cells.Item(156, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(156, 5).Value = 13
$ws.Cells.Item(156, 6).Value = 100112026
$ws.Cells.Item(156, 7).Value = "Haba"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 79
$ws.Cells.Item(156, 11).Value = 6000
$ws.Cells.Item(156, 12).Value = 7000
$ws.Cells.Item(156, 13).Value = 6494
$ws.Cells.Item(156, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(156, 15).Value = "Región Metropolitana"
$ws.Cells.Item(156, 16).Value = 260
$ws.Cells.Item(156, 17).Value = 25
$ws.Cells.Item(156, 18).Value = "Hortaliza"
